# "Update of all values to match PDF edition 10 (commit 1)"
# Appends the 2021 data row to the ATM operational-units table, updates the
# active selection, and adds the footnote defined names that accompany the
# refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data for year 2021 (columns: year, ansp, acc, app, twr, afis, atc)
$ws.Range("A12").Value = 2021
$ws.Range("B12").Value = 38
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = 276
$ws.Range("E12").Value = 381
$ws.Range("F12").Value = 81
$ws.Range("G12").Value = 639

# Footnote references (sheet-scoped defined names) tied to the new edition
$ws.Names.Add('_ftn1', '=Sheet1!#REF!')
$ws.Names.Add('_ftn2', '=Sheet1!#REF!')
$ws.Names.Add('_ftnref1', '=Sheet1!$A$6')
$ws.Names.Add('_ftnref2', '=Sheet1!$A$12')

# Leave the selection where the author last left it
[void]$ws.Range("N11").Select()
